$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.006.10'
$ws.Range('E2').Value = '  -1.47%  '
$ws.Range('D3').Value = '2.449.71'
$ws.Range('E3').Value = '  -2.15%  '
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '523.88'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.66%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '129.79'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -3.58%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.0976'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -1.76%  '
$ws.Range('E10').Value = '  -2.04%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '4.96'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -4.62%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.323'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -4.17%  '
$ws.Range('D13').Value = '2.889.70'
$ws.Range('E13').Value = '  -2.82%  '
$ws.Range('D14').Value = '57.949.77'
$ws.Range('E14').Value = '  -1.48%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '21.57'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -4.06%  '
$ws.Range('E16').Value = '  -2.95%  '
$ws.Range('D17').Value = '2.457.50'
$ws.Range('E17').Value = '  -3.18%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '10.38'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -3.57%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '4.12'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -2.34%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '314.05'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -3.39%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.16'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.34%  '
$ws.Range('E22').Value = '  +0.32%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '65.17'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.25%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.404'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -1.82%  '
$ws.Range('B25').Value = 'Binance-PegBSC-USD'
$ws.Range('C25').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.00'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +1.04%  '
$ws.Range('B26').Value = 'WrappedeETH'
$ws.Range('C26').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D26').Value = '2.568.52'
$ws.Range('E26').Value = '  -2.96%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.157'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -2.67%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '7.25'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -2.95%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '174.83'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +3.76%  '
$ws.Range('D30').Value = '0.0₃0736'
$ws.Range('E30').Value = '  -3.43%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.69'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -2.44%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '6.15'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -4.10%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.13'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -6.74%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.999'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +0.01%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.998'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +0.67%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '17.87'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -2.90%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.18'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -7.84%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.78'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -5.27%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '36.25'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -1.32%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.808'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +2.08%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.45'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -3.87%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '3.39'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -3.07%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '126.40'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -3.49%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.585'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -3.63%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '259.83'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -8.68%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '4.79'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -7.56%  '
$ws.Range('E47').Value = '  +0.25%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0492'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -3.17%  '
$ws.Range('E49').Value = '  -3.14%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '17.04'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -5.04%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '16.35'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -5.60%  '
